$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "UniformF"

# Add new row 16 data
$ws.Cells.Item(16, 1).Value = 14
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"

for ($col = 3; $col -le 13; $col++) {
    $ws.Cells.Item(16, $col).Value = 1
}
